$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{ "b" = "center"; "y" = "left"; "r" = "right" }

for ($row = 2; $row -le 361; $row++) {
    $cell = $ws.Cells.Item($row, 12)  # column L
    $val = $cell.Value()
    if ($map.ContainsKey($val)) {
        $cell.Value = $map[$val]
    }
}

for ($row = 1; $row -le 361; $row++) {
    for ($col = 1; $col -le 4; $col++) {
        $cell = $ws.Cells.Item($row, $col)
        $val = $cell.Value()
        if ($val -like "face//face_*") {
            $cell.Value = $val -replace "^face//face_", "book//book_"
        }
    }
}
